$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new date
$ws.Name = "SPY 2023-06-29"

# Remove the trailing 8 rows that no longer exist in the updated data (rows 38-45)
$ws.Rows("38:45").Delete()

# Update data rows 2-37 with refreshed option-chain data for 2023-06-29 expiry
# Row 2: SPY230629C00380000
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "SPY230629C00380000"
$ws.Cells.Item(2, 3).Value = "N/A"
$ws.Cells.Item(2, 4).Value = 380
$ws.Cells.Item(2, 5).Value = 55.9
$ws.Cells.Item(2, 6).Value = 55.29
$ws.Cells.Item(2, 7).Value = 55.54
$ws.Cells.Item(2, 8).Value = 2.3199997
$ws.Cells.Item(2, 9).Value = 4.329973
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(2, 11).Value = 16
$ws.Cells.Item(2, 12).Value = 0.00001
$ws.Cells.Item(2, 13).Value = $true
$ws.Cells.Item(2, 14).Value = "REGULAR"
$ws.Cells.Item(2, 15).Value = "USD"
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0.01
$ws.Cells.Item(2, 20).Value = -0.055

# Row 3: SPY230629C00390000
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "SPY230629C00390000"
$ws.Cells.Item(3, 3).Value = "N/A"
$ws.Cells.Item(3, 4).Value = 390
$ws.Cells.Item(3, 5).Value = 44.91
$ws.Cells.Item(3, 6).Value = 45.3
$ws.Cells.Item(3, 7).Value = 45.55
$ws.Cells.Item(3, 8).Value = -1.7999992
$ws.Cells.Item(3, 9).Value = -3.853563
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 0.00001
$ws.Cells.Item(3, 13).Value = $true
$ws.Cells.Item(3, 14).Value = "REGULAR"
$ws.Cells.Item(3, 15).Value = "USD"
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0.011
$ws.Cells.Item(3, 20).Value = -0.056

# Row 4: SPY230629C00400000
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "SPY230629C00400000"
$ws.Cells.Item(4, 3).Value = "N/A"
$ws.Cells.Item(4, 4).Value = 400
$ws.Cells.Item(4, 5).Value = 35.67
$ws.Cells.Item(4, 6).Value = 35.26
$ws.Cells.Item(4, 7).Value = 35.55
$ws.Cells.Item(4, 8).Value = 3.709999
$ws.Cells.Item(4, 9).Value = 11.608258
$ws.Cells.Item(4, 10).Value = 119
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.00001
$ws.Cells.Item(4, 13).Value = $true
$ws.Cells.Item(4, 14).Value = "REGULAR"
$ws.Cells.Item(4, 15).Value = "USD"
$ws.Cells.Item(4, 16).Value = 1
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0.011
$ws.Cells.Item(4, 20).Value = -0.058

# Row 5: SPY230629C00410000
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "SPY230629C00410000"
$ws.Cells.Item(5, 3).Value = "N/A"
$ws.Cells.Item(5, 4).Value = 410
$ws.Cells.Item(5, 5).Value = 25.59
$ws.Cells.Item(5, 6).Value = 25.3
$ws.Cells.Item(5, 7).Value = 25.55
$ws.Cells.Item(5, 8).Value = 2.8099995
$ws.Cells.Item(5, 9).Value = 12.335379
$ws.Cells.Item(5, 10).Value = 12
$ws.Cells.Item(5, 11).Value = 17
$ws.Cells.Item(5, 12).Value = 0.00001
$ws.Cells.Item(5, 13).Value = $true
$ws.Cells.Item(5, 14).Value = "REGULAR"
$ws.Cells.Item(5, 15).Value = "USD"
$ws.Cells.Item(5, 16).Value = 1
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0.011
$ws.Cells.Item(5, 20).Value = -0.059

# Row 6: SPY230629C00420000
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "SPY230629C00420000"
$ws.Cells.Item(6, 3).Value = "N/A"
$ws.Cells.Item(6, 4).Value = 420
$ws.Cells.Item(6, 5).Value = 15.85
$ws.Cells.Item(6, 6).Value = 15.39
$ws.Cells.Item(6, 7).Value = 15.65
$ws.Cells.Item(6, 8).Value = -0.69000053
$ws.Cells.Item(6, 9).Value = -4.171708
$ws.Cells.Item(6, 10).Value = 155
$ws.Cells.Item(6, 11).Value = 66
$ws.Cells.Item(6, 12).Value = 0.00001
$ws.Cells.Item(6, 13).Value = $true
$ws.Cells.Item(6, 14).Value = "REGULAR"
$ws.Cells.Item(6, 15).Value = "USD"
$ws.Cells.Item(6, 16).Value = 1
$ws.Cells.Item(6, 17).Value = 0
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0.012
$ws.Cells.Item(6, 20).Value = -0.06

# Row 7: SPY230629C00424000
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "SPY230629C00424000"
$ws.Cells.Item(7, 3).Value = "N/A"
$ws.Cells.Item(7, 4).Value = 424
$ws.Cells.Item(7, 5).Value = 11.75
$ws.Cells.Item(7, 6).Value = 11.34
$ws.Cells.Item(7, 7).Value = 11.56
$ws.Cells.Item(7, 8).Value = -1.0299997
$ws.Cells.Item(7, 9).Value = -8.059465
$ws.Cells.Item(7, 10).Value = 155
$ws.Cells.Item(7, 11).Value = 49
$ws.Cells.Item(7, 12).Value = 0.00001
$ws.Cells.Item(7, 13).Value = $true
$ws.Cells.Item(7, 14).Value = "REGULAR"
$ws.Cells.Item(7, 15).Value = "USD"
$ws.Cells.Item(7, 16).Value = 1
$ws.Cells.Item(7, 17).Value = 0
$ws.Cells.Item(7, 18).Value = 0
$ws.Cells.Item(7, 19).Value = 0.012
$ws.Cells.Item(7, 20).Value = -0.061

# Row 8: SPY230629C00426000
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "SPY230629C00426000"
$ws.Cells.Item(8, 3).Value = "N/A"
$ws.Cells.Item(8, 4).Value = 426
$ws.Cells.Item(8, 5).Value = 9.54
$ws.Cells.Item(8, 6).Value = 9.42
$ws.Cells.Item(8, 7).Value = 9.65
$ws.Cells.Item(8, 8).Value = -0.75
$ws.Cells.Item(8, 9).Value = -7.2886295
$ws.Cells.Item(8, 10).Value = 101
$ws.Cells.Item(8, 11).Value = 20
$ws.Cells.Item(8, 12).Value = 0.00001
$ws.Cells.Item(8, 13).Value = $true
$ws.Cells.Item(8, 14).Value = "REGULAR"
$ws.Cells.Item(8, 15).Value = "USD"
$ws.Cells.Item(8, 16).Value = 1
$ws.Cells.Item(8, 17).Value = 0
$ws.Cells.Item(8, 18).Value = 0
$ws.Cells.Item(8, 19).Value = 0.012
$ws.Cells.Item(8, 20).Value = -0.061

# Row 9: SPY230629C00428000
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "SPY230629C00428000"
$ws.Cells.Item(9, 3).Value = "N/A"
$ws.Cells.Item(9, 4).Value = 428
$ws.Cells.Item(9, 5).Value = 7.61
$ws.Cells.Item(9, 6).Value = 7.35
$ws.Cells.Item(9, 7).Value = 7.57
$ws.Cells.Item(9, 8).Value = -1.3200002
$ws.Cells.Item(9, 9).Value = -14.781636
$ws.Cells.Item(9, 10).Value = 68
$ws.Cells.Item(9, 11).Value = 66
$ws.Cells.Item(9, 12).Value = 0.00001
$ws.Cells.Item(9, 13).Value = $true
$ws.Cells.Item(9, 14).Value = "REGULAR"
$ws.Cells.Item(9, 15).Value = "USD"
$ws.Cells.Item(9, 16).Value = 1
$ws.Cells.Item(9, 17).Value = 0
$ws.Cells.Item(9, 18).Value = 0
$ws.Cells.Item(9, 19).Value = 0.012
$ws.Cells.Item(9, 20).Value = -0.062

# Row 10: SPY230629C00429000
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "SPY230629C00429000"
$ws.Cells.Item(10, 3).Value = "N/A"
$ws.Cells.Item(10, 4).Value = 429
$ws.Cells.Item(10, 5).Value = 6.48
$ws.Cells.Item(10, 6).Value = 6.37
$ws.Cells.Item(10, 7).Value = 6.61
$ws.Cells.Item(10, 8).Value = -0.75
$ws.Cells.Item(10, 9).Value = -10.373444
$ws.Cells.Item(10, 10).Value = 82
$ws.Cells.Item(10, 11).Value = 122
$ws.Cells.Item(10, 12).Value = 0.00001
$ws.Cells.Item(10, 13).Value = $true
$ws.Cells.Item(10, 14).Value = "REGULAR"
$ws.Cells.Item(10, 15).Value = "USD"
$ws.Cells.Item(10, 16).Value = 1
$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 0
$ws.Cells.Item(10, 19).Value = 0.012
$ws.Cells.Item(10, 20).Value = -0.062

# Row 11: SPY230629C00430000
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "SPY230629C00430000"
$ws.Cells.Item(11, 3).Value = "N/A"
$ws.Cells.Item(11, 4).Value = 430
$ws.Cells.Item(11, 5).Value = 5.59
$ws.Cells.Item(11, 6).Value = 5.39
$ws.Cells.Item(11, 7).Value = 5.6
$ws.Cells.Item(11, 8).Value = -0.71000004
$ws.Cells.Item(11, 9).Value = -11.269841
$ws.Cells.Item(11, 10).Value = 434
$ws.Cells.Item(11, 11).Value = 894
$ws.Cells.Item(11, 12).Value = 0.00001
$ws.Cells.Item(11, 13).Value = $true
$ws.Cells.Item(11, 14).Value = "REGULAR"
$ws.Cells.Item(11, 15).Value = "USD"
$ws.Cells.Item(11, 16).Value = 1
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 0.012
$ws.Cells.Item(11, 20).Value = -0.062

# Row 12: SPY230629C00431000
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "SPY230629C00431000"
$ws.Cells.Item(12, 3).Value = "N/A"
$ws.Cells.Item(12, 4).Value = 431
$ws.Cells.Item(12, 5).Value = 4.54
$ws.Cells.Item(12, 6).Value = 4.36
$ws.Cells.Item(12, 7).Value = 4.59
$ws.Cells.Item(12, 8).Value = -1.5100002
$ws.Cells.Item(12, 9).Value = -24.958681
$ws.Cells.Item(12, 10).Value = 1102
$ws.Cells.Item(12, 11).Value = 2429
$ws.Cells.Item(12, 12).Value = 0.00001
$ws.Cells.Item(12, 13).Value = $true
$ws.Cells.Item(12, 14).Value = "REGULAR"
$ws.Cells.Item(12, 15).Value = "USD"
$ws.Cells.Item(12, 16).Value = 1
$ws.Cells.Item(12, 17).Value = 0
$ws.Cells.Item(12, 18).Value = 0
$ws.Cells.Item(12, 19).Value = 0.012
$ws.Cells.Item(12, 20).Value = -0.062

# Row 13: SPY230629C00432000
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "SPY230629C00432000"
$ws.Cells.Item(13, 3).Value = "N/A"
$ws.Cells.Item(13, 4).Value = 432
$ws.Cells.Item(13, 5).Value = 3.65
$ws.Cells.Item(13, 6).Value = 3.48
$ws.Cells.Item(13, 7).Value = 3.64
$ws.Cells.Item(13, 8).Value = -0.7999997
$ws.Cells.Item(13, 9).Value = -17.977522
$ws.Cells.Item(13, 10).Value = 1185
$ws.Cells.Item(13, 11).Value = 2962
$ws.Cells.Item(13, 12).Value = 0.00001
$ws.Cells.Item(13, 13).Value = $true
$ws.Cells.Item(13, 14).Value = "REGULAR"
$ws.Cells.Item(13, 15).Value = "USD"
$ws.Cells.Item(13, 16).Value = 1
$ws.Cells.Item(13, 17).Value = 0
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 0.012
$ws.Cells.Item(13, 20).Value = -0.062

# Row 14: SPY230629C00433000
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "SPY230629C00433000"
$ws.Cells.Item(14, 3).Value = "N/A"
$ws.Cells.Item(14, 4).Value = 433
$ws.Cells.Item(14, 5).Value = 2.78
$ws.Cells.Item(14, 6).Value = 2.68
$ws.Cells.Item(14, 7).Value = 2.71
$ws.Cells.Item(14, 8).Value = -0.77
$ws.Cells.Item(14, 9).Value = -21.690142
$ws.Cells.Item(14, 10).Value = 5659
$ws.Cells.Item(14, 11).Value = 3625
$ws.Cells.Item(14, 12).Value = 0.00001
$ws.Cells.Item(14, 13).Value = $true
$ws.Cells.Item(14, 14).Value = "REGULAR"
$ws.Cells.Item(14, 15).Value = "USD"
$ws.Cells.Item(14, 16).Value = 1
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 0.012
$ws.Cells.Item(14, 20).Value = -0.062

# Row 15: SPY230629C00434000
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "SPY230629C00434000"
$ws.Cells.Item(15, 3).Value = "N/A"
$ws.Cells.Item(15, 4).Value = 434
$ws.Cells.Item(15, 5).Value = 2.03
$ws.Cells.Item(15, 6).Value = 1.97
$ws.Cells.Item(15, 7).Value = 1.99
$ws.Cells.Item(15, 8).Value = -0.77
$ws.Cells.Item(15, 9).Value = -27.5
$ws.Cells.Item(15, 10).Value = 24408
$ws.Cells.Item(15, 11).Value = 4854
$ws.Cells.Item(15, 12).Value = 0.06714800048828126
$ws.Cells.Item(15, 13).Value = $true
$ws.Cells.Item(15, 14).Value = "REGULAR"
$ws.Cells.Item(15, 15).Value = "USD"
$ws.Cells.Item(15, 16).Value = 0.878
$ws.Cells.Item(15, 17).Value = 0.132
$ws.Cells.Item(15, 18).Value = 0.046
$ws.Cells.Item(15, 19).Value = 0.01
$ws.Cells.Item(15, 20).Value = -0.21

# Row 16: SPY230629C00435000
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "SPY230629C00435000"
$ws.Cells.Item(16, 3).Value = "N/A"
$ws.Cells.Item(16, 4).Value = 435
$ws.Cells.Item(16, 5).Value = 1.42
$ws.Cells.Item(16, 6).Value = 1.39
$ws.Cells.Item(16, 7).Value = 1.4
$ws.Cells.Item(16, 8).Value = -0.72000015
$ws.Cells.Item(16, 9).Value = -33.644863
$ws.Cells.Item(16, 10).Value = 101637
$ws.Cells.Item(16, 11).Value = 7759
$ws.Cells.Item(16, 12).Value = 0.07715766601562501
$ws.Cells.Item(16, 13).Value = $true
$ws.Cells.Item(16, 14).Value = "REGULAR"
$ws.Cells.Item(16, 15).Value = "USD"
$ws.Cells.Item(16, 16).Value = 0.672
$ws.Cells.Item(16, 17).Value = 0.205
$ws.Cells.Item(16, 18).Value = 0.082
$ws.Cells.Item(16, 19).Value = 0.008
$ws.Cells.Item(16, 20).Value = -0.36

# Row 17: SPY230629C00436000
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "SPY230629C00436000"
$ws.Cells.Item(17, 3).Value = "N/A"
$ws.Cells.Item(17, 4).Value = 436
$ws.Cells.Item(17, 5).Value = 0.91
$ws.Cells.Item(17, 6).Value = 0.89
$ws.Cells.Item(17, 7).Value = 0.9
$ws.Cells.Item(17, 8).Value = -0.6499999
$ws.Cells.Item(17, 9).Value = -41.666664
$ws.Cells.Item(17, 10).Value = 199597
$ws.Cells.Item(17, 11).Value = 7028
$ws.Cells.Item(17, 12).Value = 0.07996525512695313
$ws.Cells.Item(17, 13).Value = $false
$ws.Cells.Item(17, 14).Value = "REGULAR"
$ws.Cells.Item(17, 15).Value = "USD"
$ws.Cells.Item(17, 16).Value = 0.452
$ws.Cells.Item(17, 17).Value = 0.217
$ws.Cells.Item(17, 18).Value = 0.09
$ws.Cells.Item(17, 19).Value = 0.005
$ws.Cells.Item(17, 20).Value = -0.389

# Row 18: SPY230629C00437000
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "SPY230629C00437000"
$ws.Cells.Item(18, 3).Value = "N/A"
$ws.Cells.Item(18, 4).Value = 437
$ws.Cells.Item(18, 5).Value = 0.54
$ws.Cells.Item(18, 6).Value = 0.52
$ws.Cells.Item(18, 7).Value = 0.53
$ws.Cells.Item(18, 8).Value = -0.57
$ws.Cells.Item(18, 9).Value = -51.35135
$ws.Cells.Item(18, 10).Value = 134524
$ws.Cells.Item(18, 11).Value = 4017
$ws.Cells.Item(18, 12).Value = 0.081063876953125
$ws.Cells.Item(18, 13).Value = $false
$ws.Cells.Item(18, 14).Value = "REGULAR"
$ws.Cells.Item(18, 15).Value = "USD"
$ws.Cells.Item(18, 16).Value = 0.255
$ws.Cells.Item(18, 17).Value = 0.174
$ws.Cells.Item(18, 18).Value = 0.073
$ws.Cells.Item(18, 19).Value = 0.003
$ws.Cells.Item(18, 20).Value = -0.313

# Row 19: SPY230629C00438000
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "SPY230629C00438000"
$ws.Cells.Item(19, 3).Value = "N/A"
$ws.Cells.Item(19, 4).Value = 438
$ws.Cells.Item(19, 5).Value = 0.29
$ws.Cells.Item(19, 6).Value = 0.29
$ws.Cells.Item(19, 7).Value = 0.3
$ws.Cells.Item(19, 8).Value = -0.45000002
$ws.Cells.Item(19, 9).Value = -60.81081
$ws.Cells.Item(19, 10).Value = 63264
$ws.Cells.Item(19, 11).Value = 4199
$ws.Cells.Item(19, 12).Value = 0.08326112060546878
$ws.Cells.Item(19, 13).Value = $false
$ws.Cells.Item(19, 14).Value = "REGULAR"
$ws.Cells.Item(19, 15).Value = "USD"
$ws.Cells.Item(19, 16).Value = 0.122
$ws.Cells.Item(19, 17).Value = 0.107
$ws.Cells.Item(19, 18).Value = 0.046
$ws.Cells.Item(19, 19).Value = 0.001
$ws.Cells.Item(19, 20).Value = -0.2

# Row 20: SPY230629C00439000
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "SPY230629C00439000"
$ws.Cells.Item(20, 3).Value = "N/A"
$ws.Cells.Item(20, 4).Value = 439
$ws.Cells.Item(20, 5).Value = 0.16
$ws.Cells.Item(20, 6).Value = 0.14
$ws.Cells.Item(20, 7).Value = 0.15
$ws.Cells.Item(20, 8).Value = -0.34
$ws.Cells.Item(20, 9).Value = -68
$ws.Cells.Item(20, 10).Value = 51433
$ws.Cells.Item(20, 11).Value = 2929
$ws.Cells.Item(20, 12).Value = 0.08374939697265626
$ws.Cells.Item(20, 13).Value = $false
$ws.Cells.Item(20, 14).Value = "REGULAR"
$ws.Cells.Item(20, 15).Value = "USD"
$ws.Cells.Item(20, 16).Value = 0.047
$ws.Cells.Item(20, 17).Value = 0.051
$ws.Cells.Item(20, 18).Value = 0.022
$ws.Cells.Item(20, 19).Value = 0.001
$ws.Cells.Item(20, 20).Value = -0.096

# Row 21: SPY230629C00440000
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "SPY230629C00440000"
$ws.Cells.Item(21, 3).Value = "N/A"
$ws.Cells.Item(21, 4).Value = 440
$ws.Cells.Item(21, 5).Value = 0.08
$ws.Cells.Item(21, 6).Value = 0.07
$ws.Cells.Item(21, 7).Value = 0.08
$ws.Cells.Item(21, 8).Value = -0.23
$ws.Cells.Item(21, 9).Value = -74.19355
$ws.Cells.Item(21, 10).Value = 75854
$ws.Cells.Item(21, 11).Value = 8585
$ws.Cells.Item(21, 12).Value = 0.08692319335937501
$ws.Cells.Item(21, 13).Value = $false
$ws.Cells.Item(21, 14).Value = "REGULAR"
$ws.Cells.Item(21, 15).Value = "USD"
$ws.Cells.Item(21, 16).Value = 0.017
$ws.Cells.Item(21, 17).Value = 0.021
$ws.Cells.Item(21, 18).Value = 0.01
$ws.Cells.Item(21, 19).Value = 0
$ws.Cells.Item(21, 20).Value = -0.043

# Row 22: SPY230629C00441000
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "SPY230629C00441000"
$ws.Cells.Item(22, 3).Value = "N/A"
$ws.Cells.Item(22, 4).Value = 441
$ws.Cells.Item(22, 5).Value = 0.05
$ws.Cells.Item(22, 6).Value = 0.04
$ws.Cells.Item(22, 7).Value = 0.05
$ws.Cells.Item(22, 8).Value = -0.13000001
$ws.Cells.Item(22, 9).Value = -72.22223
$ws.Cells.Item(22, 10).Value = 23319
$ws.Cells.Item(22, 11).Value = 2061
$ws.Cells.Item(22, 12).Value = 0.0932707861328125
$ws.Cells.Item(22, 13).Value = $false
$ws.Cells.Item(22, 14).Value = "REGULAR"
$ws.Cells.Item(22, 15).Value = "USD"
$ws.Cells.Item(22, 16).Value = 0.007
$ws.Cells.Item(22, 17).Value = 0.01
$ws.Cells.Item(22, 18).Value = 0.005
$ws.Cells.Item(22, 19).Value = 0
$ws.Cells.Item(22, 20).Value = -0.022

# Row 23: SPY230629C00442000
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "SPY230629C00442000"
$ws.Cells.Item(23, 3).Value = "N/A"
$ws.Cells.Item(23, 4).Value = 442
$ws.Cells.Item(23, 5).Value = 0.02
$ws.Cells.Item(23, 6).Value = 0.02
$ws.Cells.Item(23, 7).Value = 0.03
$ws.Cells.Item(23, 8).Value = -0.09
$ws.Cells.Item(23, 9).Value = -81.818184
$ws.Cells.Item(23, 10).Value = 9819
$ws.Cells.Item(23, 11).Value = 3510
$ws.Cells.Item(23, 12).Value = 0.098641826171875
$ws.Cells.Item(23, 13).Value = $false
$ws.Cells.Item(23, 14).Value = "REGULAR"
$ws.Cells.Item(23, 15).Value = "USD"
$ws.Cells.Item(23, 16).Value = 0.003
$ws.Cells.Item(23, 17).Value = 0.004
$ws.Cells.Item(23, 18).Value = 0.002
$ws.Cells.Item(23, 19).Value = 0
$ws.Cells.Item(23, 20).Value = -0.011

# Row 24: SPY230629C00443000
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "SPY230629C00443000"
$ws.Cells.Item(24, 3).Value = "N/A"
$ws.Cells.Item(24, 4).Value = 443
$ws.Cells.Item(24, 5).Value = 0.02
$ws.Cells.Item(24, 6).Value = 0.01
$ws.Cells.Item(24, 7).Value = 0.02
$ws.Cells.Item(24, 8).Value = -0.04
$ws.Cells.Item(24, 9).Value = -66.66667
$ws.Cells.Item(24, 10).Value = 3613
$ws.Cells.Item(24, 11).Value = 3059
$ws.Cells.Item(24, 12).Value = 0.104501142578125
$ws.Cells.Item(24, 13).Value = $false
$ws.Cells.Item(24, 14).Value = "REGULAR"
$ws.Cells.Item(24, 15).Value = "USD"
$ws.Cells.Item(24, 16).Value = 0.001
$ws.Cells.Item(24, 17).Value = 0.002
$ws.Cells.Item(24, 18).Value = 0.001
$ws.Cells.Item(24, 19).Value = 0
$ws.Cells.Item(24, 20).Value = -0.005

# Row 25: SPY230629C00444000
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "SPY230629C00444000"
$ws.Cells.Item(25, 3).Value = "N/A"
$ws.Cells.Item(25, 4).Value = 444
$ws.Cells.Item(25, 5).Value = 0.02
$ws.Cells.Item(25, 6).Value = 0.01
$ws.Cells.Item(25, 7).Value = 0.02
$ws.Cells.Item(25, 8).Value = -0.02
$ws.Cells.Item(25, 9).Value = -50
$ws.Cells.Item(25, 10).Value = 4028
$ws.Cells.Item(25, 11).Value = 2043
$ws.Cells.Item(25, 12).Value = 0.117196328125
$ws.Cells.Item(25, 13).Value = $false
$ws.Cells.Item(25, 14).Value = "REGULAR"
$ws.Cells.Item(25, 15).Value = "USD"
$ws.Cells.Item(25, 16).Value = 0.001
$ws.Cells.Item(25, 17).Value = 0.001
$ws.Cells.Item(25, 18).Value = 0.001
$ws.Cells.Item(25, 19).Value = 0
$ws.Cells.Item(25, 20).Value = -0.005

# Row 26: SPY230629C00445000
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "SPY230629C00445000"
$ws.Cells.Item(26, 3).Value = "N/A"
$ws.Cells.Item(26, 4).Value = 445
$ws.Cells.Item(26, 5).Value = 0.02
$ws.Cells.Item(26, 6).Value = 0.01
$ws.Cells.Item(26, 7).Value = 0.02
$ws.Cells.Item(26, 8).Value = -0.01
$ws.Cells.Item(26, 9).Value = -33.333336
$ws.Cells.Item(26, 10).Value = 1490
$ws.Cells.Item(26, 11).Value = 2163
$ws.Cells.Item(26, 12).Value = 0.1289149609375
$ws.Cells.Item(26, 13).Value = $false
$ws.Cells.Item(26, 14).Value = "REGULAR"
$ws.Cells.Item(26, 15).Value = "USD"
$ws.Cells.Item(26, 16).Value = 0.001
$ws.Cells.Item(26, 17).Value = 0.001
$ws.Cells.Item(26, 18).Value = 0.001
$ws.Cells.Item(26, 19).Value = 0
$ws.Cells.Item(26, 20).Value = -0.005

# Row 27: SPY230629C00446000
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "SPY230629C00446000"
$ws.Cells.Item(27, 3).Value = "N/A"
$ws.Cells.Item(27, 4).Value = 446
$ws.Cells.Item(27, 5).Value = 0.01
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 0.01
$ws.Cells.Item(27, 8).Value = -0.01
$ws.Cells.Item(27, 9).Value = -50
$ws.Cells.Item(27, 10).Value = 543
$ws.Cells.Item(27, 11).Value = 1145
$ws.Cells.Item(27, 12).Value = 0.1289149609375
$ws.Cells.Item(27, 13).Value = $false
$ws.Cells.Item(27, 14).Value = "REGULAR"
$ws.Cells.Item(27, 15).Value = "USD"
$ws.Cells.Item(27, 16).Value = 0
$ws.Cells.Item(27, 17).Value = 0
$ws.Cells.Item(27, 18).Value = 0
$ws.Cells.Item(27, 19).Value = 0
$ws.Cells.Item(27, 20).Value = -0.002

# Row 28: SPY230629C00447000
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "SPY230629C00447000"
$ws.Cells.Item(28, 3).Value = "N/A"
$ws.Cells.Item(28, 4).Value = 447
$ws.Cells.Item(28, 5).Value = 0.01
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0.01
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 3514
$ws.Cells.Item(28, 11).Value = 2761
$ws.Cells.Item(28, 12).Value = 0.14063359375
$ws.Cells.Item(28, 13).Value = $false
$ws.Cells.Item(28, 14).Value = "REGULAR"
$ws.Cells.Item(28, 15).Value = "USD"
$ws.Cells.Item(28, 16).Value = 0
$ws.Cells.Item(28, 17).Value = 0
$ws.Cells.Item(28, 18).Value = 0
$ws.Cells.Item(28, 19).Value = 0
$ws.Cells.Item(28, 20).Value = -0.002

# Row 29: SPY230629C00448000
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "SPY230629C00448000"
$ws.Cells.Item(29, 3).Value = "N/A"
$ws.Cells.Item(29, 4).Value = 448
$ws.Cells.Item(29, 5).Value = 0.01
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 0.01
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 1467
$ws.Cells.Item(29, 11).Value = 1308
$ws.Cells.Item(29, 12).Value = 0.15039912109375
$ws.Cells.Item(29, 13).Value = $false
$ws.Cells.Item(29, 14).Value = "REGULAR"
$ws.Cells.Item(29, 15).Value = "USD"
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0
$ws.Cells.Item(29, 18).Value = 0
$ws.Cells.Item(29, 19).Value = 0
$ws.Cells.Item(29, 20).Value = -0.001

# Row 30: SPY230629C00450000
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "SPY230629C00450000"
$ws.Cells.Item(30, 3).Value = "N/A"
$ws.Cells.Item(30, 4).Value = 450
$ws.Cells.Item(30, 5).Value = 0.01
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0.01
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 68
$ws.Cells.Item(30, 11).Value = 2122
$ws.Cells.Item(30, 12).Value = 0.17188328125
$ws.Cells.Item(30, 13).Value = $false
$ws.Cells.Item(30, 14).Value = "REGULAR"
$ws.Cells.Item(30, 15).Value = "USD"
$ws.Cells.Item(30, 16).Value = 0
$ws.Cells.Item(30, 17).Value = 0
$ws.Cells.Item(30, 18).Value = 0
$ws.Cells.Item(30, 19).Value = 0
$ws.Cells.Item(30, 20).Value = -0.001

# Row 31: SPY230629C00452000
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "SPY230629C00452000"
$ws.Cells.Item(31, 3).Value = "N/A"
$ws.Cells.Item(31, 4).Value = 452
$ws.Cells.Item(31, 5).Value = 0.01
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0.01
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 4
$ws.Cells.Item(31, 11).Value = 20
$ws.Cells.Item(31, 12).Value = 0.1914143359375
$ws.Cells.Item(31, 13).Value = $false
$ws.Cells.Item(31, 14).Value = "REGULAR"
$ws.Cells.Item(31, 15).Value = "USD"
$ws.Cells.Item(31, 16).Value = 0
$ws.Cells.Item(31, 17).Value = 0
$ws.Cells.Item(31, 18).Value = 0
$ws.Cells.Item(31, 19).Value = 0
$ws.Cells.Item(31, 20).Value = -0.001

# Row 32: SPY230629C00454000
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = "SPY230629C00454000"
$ws.Cells.Item(32, 3).Value = "N/A"
$ws.Cells.Item(32, 4).Value = 454
$ws.Cells.Item(32, 5).Value = 0.01
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 0.01
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 18
$ws.Cells.Item(32, 11).Value = 2021
$ws.Cells.Item(32, 12).Value = 0.210945390625
$ws.Cells.Item(32, 13).Value = $false
$ws.Cells.Item(32, 14).Value = "REGULAR"
$ws.Cells.Item(32, 15).Value = "USD"
$ws.Cells.Item(32, 16).Value = 0
$ws.Cells.Item(32, 17).Value = 0
$ws.Cells.Item(32, 18).Value = 0
$ws.Cells.Item(32, 19).Value = 0
$ws.Cells.Item(32, 20).Value = -0.001

# Row 33: SPY230629C00456000
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = "SPY230629C00456000"
$ws.Cells.Item(33, 3).Value = "N/A"
$ws.Cells.Item(33, 4).Value = 456
$ws.Cells.Item(33, 5).Value = 0.01
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 0.01
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 5
$ws.Cells.Item(33, 11).Value = 295
$ws.Cells.Item(33, 12).Value = 0.23438265625
$ws.Cells.Item(33, 13).Value = $false
$ws.Cells.Item(33, 14).Value = "REGULAR"
$ws.Cells.Item(33, 15).Value = "USD"
$ws.Cells.Item(33, 16).Value = 0
$ws.Cells.Item(33, 17).Value = 0
$ws.Cells.Item(33, 18).Value = 0
$ws.Cells.Item(33, 19).Value = 0
$ws.Cells.Item(33, 20).Value = -0.001

# Row 34: SPY230629C00458000
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = "SPY230629C00458000"
$ws.Cells.Item(34, 3).Value = "N/A"
$ws.Cells.Item(34, 4).Value = 458
$ws.Cells.Item(34, 5).Value = 0.01
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0.01
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 1
$ws.Cells.Item(34, 11).Value = 31
$ws.Cells.Item(34, 12).Value = 0.2539137109375
$ws.Cells.Item(34, 13).Value = $false
$ws.Cells.Item(34, 14).Value = "REGULAR"
$ws.Cells.Item(34, 15).Value = "USD"
$ws.Cells.Item(34, 16).Value = 0
$ws.Cells.Item(34, 17).Value = 0
$ws.Cells.Item(34, 18).Value = 0
$ws.Cells.Item(34, 19).Value = 0
$ws.Cells.Item(34, 20).Value = -0.001

# Row 35: SPY230629C00460000
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = "SPY230629C00460000"
$ws.Cells.Item(35, 3).Value = "N/A"
$ws.Cells.Item(35, 4).Value = 460
$ws.Cells.Item(35, 5).Value = 0.01
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 0.01
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 1527
$ws.Cells.Item(35, 12).Value = 0.273444765625
$ws.Cells.Item(35, 13).Value = $false
$ws.Cells.Item(35, 14).Value = "REGULAR"
$ws.Cells.Item(35, 15).Value = "USD"
$ws.Cells.Item(35, 16).Value = 0
$ws.Cells.Item(35, 17).Value = 0
$ws.Cells.Item(35, 18).Value = 0
$ws.Cells.Item(35, 19).Value = 0
$ws.Cells.Item(35, 20).Value = -0.001

# Row 36: SPY230629C00470000
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = "SPY230629C00470000"
$ws.Cells.Item(36, 3).Value = "N/A"
$ws.Cells.Item(36, 4).Value = 470
$ws.Cells.Item(36, 5).Value = 0.01
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 0.01
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 4
$ws.Cells.Item(36, 11).Value = 220
$ws.Cells.Item(36, 12).Value = 0.367193828125
$ws.Cells.Item(36, 13).Value = $false
$ws.Cells.Item(36, 14).Value = "REGULAR"
$ws.Cells.Item(36, 15).Value = "USD"
$ws.Cells.Item(36, 16).Value = 0
$ws.Cells.Item(36, 17).Value = 0
$ws.Cells.Item(36, 18).Value = 0
$ws.Cells.Item(36, 19).Value = 0
$ws.Cells.Item(36, 20).Value = -0.001

# Row 37: SPY230629C00480000
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = "SPY230629C00480000"
$ws.Cells.Item(37, 3).Value = "N/A"
$ws.Cells.Item(37, 4).Value = 480
$ws.Cells.Item(37, 5).Value = 0.01
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0.01
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 1
$ws.Cells.Item(37, 11).Value = 255
$ws.Cells.Item(37, 12).Value = 0.460942890625
$ws.Cells.Item(37, 13).Value = $false
$ws.Cells.Item(37, 14).Value = "REGULAR"
$ws.Cells.Item(37, 15).Value = "USD"
$ws.Cells.Item(37, 16).Value = 0
$ws.Cells.Item(37, 17).Value = 0
$ws.Cells.Item(37, 18).Value = 0
$ws.Cells.Item(37, 19).Value = 0
$ws.Cells.Item(37, 20).Value = -0.001
